# Update countries & provincias Spain
# Refreshes the COVID "Pais" dashboard: new snapshot timestamp, updated
# per-country counters, and the resulting reshuffle of Libia/Senegal
# (same totals column, just new numbers moving Libia ahead of Senegal).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Snapshot timestamp (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 31 de Agosto de 2020 a las 13:35"

# Country labels that swap position (row 86 <-> row 87) because Libia's
# updated total now outranks Senegal's.
$ws.Cells.Item(86, 1).Value = "Libia"
$ws.Cells.Item(87, 1).Value = "Senegal"

# row -> @{ column letter = new value }
$updates = @{
    23  = @{ B = 243463; C = 168;                 E = 16257 }
    44  = @{ B = 71843;  C = 156; D = 70468; E = 694;  G = 5; H = 681 }
    61  = @{ B = 42177;  C = 163;                 E = 4372 }
    64  = @{ B = 39460;  C = 899; D = 21410; E = 17822; G = 7; H = 228 }
    75  = @{ B = 22729;  C = 525; D = 15056; E = 7521 }
    84  = @{ B = 14863;  C = 20;  D = 13832; E = 839;  G = 1; H = 192 }
    86  = @{ B = 13966;  C = 543; D = 1459;  E = 12270; G = 5; H = 237 }
    87  = @{ B = 13611;  C = 55;  D = 9439;  E = 3888;        H = 284 }
    89  = @{ B = 12097;  C = 72;  D = 11469; E = 340;  G = 1; H = 288 }
    101 = @{ B = 8086;   C = 9;                   E = 550;  G = 1; H = 336 }
    146 = @{ B = 1883;   C = 21;  D = 1400;  E = 471 }
    163 = @{ B = 1044;   C = 4;   D = 707;   E = 303;  G = 2; H = 34 }
}

$colIndex = @{ B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8 }

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Cells.Item([int]$row, $colIndex[$col]).Value = $cols[$col]
    }
}
